$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-06 Friday" "2025-06-07 Saturday"

Replace-Text "956÷2=" "160÷4="
Replace-Text "659÷6=" "642÷4="
Replace-Text "758÷5=" "657÷7="
Replace-Text "555÷3=" "569÷7="
Replace-Text "783÷4=" "155÷6="
Replace-Text "860÷5=" "778÷9="
Replace-Text "108÷4=" "654÷6="
Replace-Text "350÷4=" "858÷5="
Replace-Text "724÷3=" "904÷2="
Replace-Text "749÷6=" "743÷6="
Replace-Text "657÷8=" "777÷9="
Replace-Text "669÷9=" "580÷4="
Replace-Text "327÷5=" "915÷2="
Replace-Text "969÷7=" "662÷4="
Replace-Text "437÷8=" "483÷3="
Replace-Text "146÷4=" "443÷2="
Replace-Text "183÷9=" "398÷6="
Replace-Text "219÷8=" "261÷7="
Replace-Text "222÷9=" "606÷4="
Replace-Text "170÷8=" "285÷8="
Replace-Text "343÷2=" "310÷6="
Replace-Text "429÷9=" "491÷7="
Replace-Text "236÷7=" "922÷6="
Replace-Text "842÷6=" "529÷4="
Replace-Text "132÷9=" "469÷2="
